# Update task data used in testing (inputTrain_1.7.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 value updates
$ws.Range("D3").Value = 2
$ws.Range("F3").Value = -3
$ws.Range("H3").Value = 46

# Move the active selection from D5 to H5
$ws.Range("H5").Select()
